# Update countries & provincias Spain
# - Refresh "last updated" timestamp
# - Update Austria's daily figures
# - Insert fresh Finlandia data (pushing Indonesia/Sudafrica down a row,
#   their own totals unchanged) ahead of Indonesia/Sudafrica in the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / timestamp cell
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 10:20"

# Austria (row 17) - updated totals
$ws.Range("B17").Value = 8958
$ws.Range("C17").Value = 170
$ws.Range("E17").Value = 8393

# Row 40: Finlandia (new data, now listed first of this trio)
$ws.Range("A40").Value = "Finlandia"
$ws.Range("B40").Value = 1286
$ws.Range("C40").Value = 46
$ws.Range("D40").Value = 10
$ws.Range("E40").Value = 1265
$ws.Range("F40").Value = 32
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 11

# Row 41: Indonesia (values unchanged, shifted down one row)
$ws.Range("A41").Value = "Indonesia"
$ws.Range("B41").Value = 1285
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 64
$ws.Range("E41").Value = 1107
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 114

# Row 42: Sudafrica (values unchanged, shifted down one row)
$ws.Range("A42").Value = "Sudafrica"
$ws.Range("B42").Value = 1280
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 31
$ws.Range("E42").Value = 1247
$ws.Range("F42").Value = 7
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 2
